{"js": "// Replace the date line and the 25 \"A\u00f7B=C, D\" division answers with new\n// values, in document order. The table shape (rows/cells) is unchanged \u2014\n// only the text inside the existing paragraphs/runs is updated, so we can\n// safely walk body.paragraphs in order and swap the text of every\n// non-empty paragraph for its replacement, preserving run formatting via\n// insertText(..., Replace).\nconst replacements = [\n  \"2023-12-23 Saturday\",\n  \"63\u00f72=31, 1\",\n  \"87\u00f78=10, 7\",\n  \"43\u00f75=8, 3\",\n  \"60\u00f76=10, 0\",\n  \"62\u00f78=7, 6\",\n  \"60\u00f75=12, 0\",\n  \"96\u00f75=19, 1\",\n  \"68\u00f72=34, 0\",\n  \"81\u00f76=13, 3\",\n  \"61\u00f78=7, 5\",\n  \"90\u00f76=15, 0\",\n  \"67\u00f77=9, 4\",\n  \"32\u00f79=3, 5\",\n  \"62\u00f78=7, 6\",\n  \"75\u00f78=9, 3\",\n  \"14\u00f75=2, 4\",\n  \"43\u00f78=5, 3\",\n  \"53\u00f76=8, 5\",\n  \"52\u00f72=26, 0\",\n  \"78\u00f77=11, 1\",\n  \"47\u00f78=5, 7\",\n  \"89\u00f73=29, 2\",\n  \"35\u00f73=11, 2\",\n  \"81\u00f76=13, 3\",\n  \"76\u00f74=19, 0\",\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet i = 0;\nfor (const p of paragraphs.items) {\n  if (p.text !== \"\") {\n    if (i >= replacements.length) {\n      break;\n    }\n    p.insertText(replacements[i], Word.InsertLocation.replace);\n    i++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 25 \"A\u00f7B=C, D\" division answers in the\n# table with new values. The table shape (20 rows x 5 cols, data only in\n# rows 1/5/9/13/17) is unchanged, so every new value is written by its\n# fixed (row, column) position - this avoids any ambiguity from values\n# that are reused elsewhere in the grid (e.g. \"81\u00f76=13, 3\" is both an old\n# and a new value at two different cells).\n$d = $word.ActiveDocument\n\n$d.Paragraphs(1).Range.Text = \"2023-12-23 Saturday\"\n\n$tbl = $d.Tables(1)\n\n$values = @{\n    \"1,1\" = \"63\u00f72=31, 1\"\n    \"1,2\" = \"87\u00f78=10, 7\"\n    \"1,3\" = \"43\u00f75=8, 3\"\n    \"1,4\" = \"60\u00f76=10, 0\"\n    \"1,5\" = \"62\u00f78=7, 6\"\n    \"5,1\" = \"60\u00f75=12, 0\"\n    \"5,2\" = \"96\u00f75=19, 1\"\n    \"5,3\" = \"68\u00f72=34, 0\"\n    \"5,4\" = \"81\u00f76=13, 3\"\n    \"5,5\" = \"61\u00f78=7, 5\"\n    \"9,1\" = \"90\u00f76=15, 0\"\n    \"9,2\" = \"67\u00f77=9, 4\"\n    \"9,3\" = \"32\u00f79=3, 5\"\n    \"9,4\" = \"62\u00f78=7, 6\"\n    \"9,5\" = \"75\u00f78=9, 3\"\n    \"13,1\" = \"14\u00f75=2, 4\"\n    \"13,2\" = \"43\u00f78=5, 3\"\n    \"13,3\" = \"53\u00f76=8, 5\"\n    \"13,4\" = \"52\u00f72=26, 0\"\n    \"13,5\" = \"78\u00f77=11, 1\"\n    \"17,1\" = \"47\u00f78=5, 7\"\n    \"17,2\" = \"89\u00f73=29, 2\"\n    \"17,3\" = \"35\u00f73=11, 2\"\n    \"17,4\" = \"81\u00f76=13, 3\"\n    \"17,5\" = \"76\u00f74=19, 0\"\n}\n\n$dataRows = @(1, 5, 9, 13, 17)\nforeach ($r in $dataRows) {\n    for ($c = 1; $c -le 5; $c++) {\n        $key = \"$r,$c\"\n        $cell = $tbl.Cell($r, $c)\n        $cell.Range.Text = $values[$key]\n    }\n}\n"}
